$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.850993812084198
$ws.Range("B1").Value = 1.285834670066833
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.360522508621216
